$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value: B3 changes from 5500 to 6000
$ws.Range("B3").Value = 6000

# Add new rows 6-8
$ws.Range("A6").Value = "Tv5 C4"
$ws.Range("B6").Value = 3200
$ws.Range("C6").Value = 313

$ws.Range("A7").Value = "МЕФ Мяу"
$ws.Range("B7").Value = 2400
$ws.Range("C7").Value = 277

$ws.Range("A8").Value = "МЕФ Мяу"
$ws.Range("B8").Value = 5000
$ws.Range("C8").Value = 278

# Update selection to C8 to mirror the final cursor position
$ws.Range("C8").Select()
